# Weekly price update: insert the newest week's record as a new row 75,
# pushing the existing rows 75-102 down to 76-103.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 75 - this shifts rows 75..102 down to 76..103
# and grows the used range from A1:R102 to A1:R103 automatically.
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with this week's data.
$ws.Cells.Item(75, 1).Value = 2
$ws.Cells.Item(75, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(75, 3).Value = "Coquimbo"
$ws.Cells.Item(75, 4).Value = 44559
$ws.Cells.Item(75, 5).Value = 4
$ws.Cells.Item(75, 6).Value = 100112024
$ws.Cells.Item(75, 7).Value = "Choclo"
$ws.Cells.Item(75, 8).Value = "Dulce o Americano"
$ws.Cells.Item(75, 9).Value = "Primera"
$ws.Cells.Item(75, 10).Value = 600
$ws.Cells.Item(75, 11).Value = 11000
$ws.Cells.Item(75, 12).Value = 12000
$ws.Cells.Item(75, 13).Value = 11500
$ws.Cells.Item(75, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(75, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(75, 16).Value = 164
$ws.Cells.Item(75, 17).Value = 70
$ws.Cells.Item(75, 18).Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Cells.Item(75, 4).NumberFormat = $ws.Cells.Item(76, 4).NumberFormat
